$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text-valued cells (Coin names and Link URLs) - safe to assign directly
$ws.Range("B20").Value = 'InternetComputer(DFINITY)'
$ws.Range("C20").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'

# Price cells (column D) - force Text number format first so Excel does not
# reinterpret numeric-looking strings (with single decimal points) as numbers,
# which would drop formatting like trailing zeros.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '39.688.90'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.205.62'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '291.85'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '86.31'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.515'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.470'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '30.27'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0784'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.42'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.543.40'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.00'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.198.54'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.727'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '39.635.36'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.35'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0877'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.78'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.60'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '235.50'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.46'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.64'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '151.70'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '32.50'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.91'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.78'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '15.96'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0985'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.79'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.063.67'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.17'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.93'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '17.73'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.60'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.418.73'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '70.74'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '88.72'

# Volume(1h) percentage cells (column E) - these contain surrounding spaces
# and a percent sign so Excel keeps them as text automatically.
$ws.Range("E2").Value = '  +1.34%  '
$ws.Range("E3").Value = '  +0.30%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("E5").Value = '  -1.10%  '
$ws.Range("E6").Value = '  +6.64%  '
$ws.Range("E7").Value = '  +1.12%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  +0.65%  '
$ws.Range("E10").Value = '  +4.34%  '
$ws.Range("E11").Value = '  +2.27%  '
$ws.Range("E12").Value = '  +1.28%  '
$ws.Range("E13").Value = '  +1.58%  '
$ws.Range("E14").Value = '  +1.78%  '
$ws.Range("E15").Value = '  +0.08%  '
$ws.Range("E16").Value = '  +0.61%  '
$ws.Range("E17").Value = '  +0.46%  '
$ws.Range("E18").Value = '  +2.49%  '
$ws.Range("E19").Value = '  +1.43%  '
$ws.Range("E20").Value = '  +10.74%  '
$ws.Range("E21").Value = '  +1.21%  '
$ws.Range("E22").Value = '  +1.55%  '
$ws.Range("E23").Value = '  +1.48%  '
$ws.Range("E24").Value = '  +4.73%  '
$ws.Range("E25").Value = '  -0.07%  '
$ws.Range("E26").Value = '  +2.80%  '
$ws.Range("E27").Value = '  +2.41%  '
$ws.Range("E28").Value = '  +0.64%  '
$ws.Range("E29").Value = '  +1.03%  '
$ws.Range("E30").Value = '  +2.39%  '
$ws.Range("E31").Value = '  +2.09%  '
$ws.Range("E32").Value = '  +2.92%  '
$ws.Range("E33").Value = '  -0.04%  '
$ws.Range("E34").Value = '  +2.78%  '
$ws.Range("E35").Value = '  +3.11%  '
$ws.Range("E37").Value = '  +2.29%  '
$ws.Range("E38").Value = '  +6.36%  '
$ws.Range("E39").Value = '  +4.81%  '
$ws.Range("E40").Value = '  +3.45%  '
$ws.Range("E41").Value = '  +3.22%  '
$ws.Range("E42").Value = '  +5.49%  '
$ws.Range("E43").Value = '  +8.86%  '
$ws.Range("E44").Value = '  +4.79%  '
$ws.Range("E45").Value = '  +3.34%  '
$ws.Range("E46").Value = '  +10.13%  '
$ws.Range("E47").Value = '  +10.87%  '
$ws.Range("E48").Value = '  +0.03%  '
$ws.Range("E49").Value = '  +0.42%  '
$ws.Range("E50").Value = '  -1.00%  '
$ws.Range("E51").Value = '  +1.92%  '
